$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.09008902633495776, 0.2712929561332827, 0.1215236909347529, 0.3486024826858709, 0.3459885601116726, 19)
    3  = @(0.440854564544955, 0.4977539282968331, 0.4961366786238803, 0.7043697030848788, 0.5652750901259675, 18)
    4  = @(0.6924743815712978, 0.7117042034689857, 0.9527927106825833, 0.9761110135033737, 0.7091201857914559, 17)
    5  = @(0.6938486294701675, 0.696168599774547, 0.7176852045027668, 0.8471630330123989, 0.5020058792592104, 16)
    6  = @(0.5732877799022984, 0.5732877799022984, 0.438377994181727, 0.6621011963300829, 0.3428646894649251, 15)
    7  = @(0.4864587127141264, 0.4881921051911807, 0.3241750264404118, 0.5693637031286871, 0.3070280946439436, 14)
    8  = @(0.4121643343351353, 0.4244592671431573, 0.2374407211537076, 0.4872788946319219, 0.2705390721203952, 13)
    9  = @(0.4514641206974861, 0.4514641206974861, 0.2461149259357549, 0.4960997137025528, 0.2181328903808997, 9)
    10 = @(0.3763340750174313, 0.3763340750174313, 0.1542216358128203, 0.3927106260503022, 0.1254706130613591, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
